$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to text format so numeric-looking strings (e.g. "24.43")
# are preserved as text instead of being auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '29.328.24'
$ws.Range("E2").Value = '  -0.01%  '
$ws.Range("D3").Value = '1.844.01'
$ws.Range("E3").Value = '  -0.14%  '
$ws.Range("E4").Value = '  +0.16%  '
$ws.Range("D5").Value = '240.09'
$ws.Range("E5").Value = '  +0.03%  '
$ws.Range("E6").Value = '  +0.16%  '
$ws.Range("D7").Value = '1.001'
$ws.Range("E7").Value = '  +0.20%  '
$ws.Range("D8").Value = '0.07450'
$ws.Range("E8").Value = '  -1.91%  '
$ws.Range("E9").Value = '  -0.40%  '
$ws.Range("D10").Value = '24.43'
$ws.Range("E10").Value = '  -0.85%  '
$ws.Range("E11").Value = '  -0.01%  '
$ws.Range("D12").Value = '1.844.11'
$ws.Range("E12").Value = '  -2.34%  '
$ws.Range("D13").Value = '4.982'
$ws.Range("E13").Value = '  -0.72%  '
$ws.Range("D14").Value = '0.6788'
$ws.Range("E14").Value = '  +0.02%  '
$ws.Range("D15").Value = '0.00001041'
$ws.Range("E15").Value = '  -1.35%  '
$ws.Range("D16").Value = '81.85'
$ws.Range("E16").Value = '  -1.32%  '
$ws.Range("D17").Value = '6.184'
$ws.Range("E17").Value = '  +0.97%  '
$ws.Range("D18").Value = '29.383.48'
$ws.Range("E18").Value = '  +0.02%  '
$ws.Range("D19").Value = '227.81'
$ws.Range("E19").Value = '  -0.58%  '
$ws.Range("E20").Value = '  -0.24%  '
$ws.Range("E21").Value = '  +0.18%  '
$ws.Range("D22").Value = '7.501'
$ws.Range("E22").Value = '  +0.44%  '
$ws.Range("D23").Value = '1.001'
$ws.Range("E23").Value = '  +0.23%  '
$ws.Range("D24").Value = '159.26'
$ws.Range("E24").Value = '  +0.52%  '
$ws.Range("D25").Value = '8.467'
$ws.Range("E25").Value = '  +0.43%  '
$ws.Range("D26").Value = '0.1366'
$ws.Range("E26").Value = '  -1.38%  '
$ws.Range("E27").Value = '  -1.00%  '
$ws.Range("D28").Value = '0.06526'
$ws.Range("E28").Value = '  +16.56%  '
$ws.Range("E29").Value = '  -1.19%  '
$ws.Range("E30").Value = '  +1.32%  '
$ws.Range("D31").Value = '4.084'
$ws.Range("E31").Value = '  -0.44%  '
$ws.Range("D32").Value = '4.076'
$ws.Range("E32").Value = '  +0.34%  '
$ws.Range("D33").Value = '1.826'
$ws.Range("E33").Value = '  -0.11%  '
$ws.Range("D34").Value = '1.140'
$ws.Range("E34").Value = '  -1.68%  '
$ws.Range("D35").Value = '0.6949'
$ws.Range("E35").Value = '  -0.02%  '
$ws.Range("D36").Value = '2.582'
$ws.Range("E36").Value = '  +0.00%  '
$ws.Range("D37").Value = '1.261.60'
$ws.Range("E37").Value = '  +2.36%  '
$ws.Range("E38").Value = '  +3.97%  '
$ws.Range("D39").Value = '0.01831'
$ws.Range("E39").Value = '  +1.85%  '
$ws.Range("D40").Value = '6.784'
$ws.Range("E40").Value = '  +6.51%  '
$ws.Range("D41").Value = '0.9195'
$ws.Range("E41").Value = '  +1.80%  '
$ws.Range("D42").Value = '0.9994'
$ws.Range("E42").Value = '  +0.09%  '
$ws.Range("D43").Value = '2.004.34'
$ws.Range("E43").Value = '  +1.29%  '
$ws.Range("D44").Value = '101.51'
$ws.Range("D45").Value = '65.93'
$ws.Range("B46").Value = 'BabyDogeCoin'
$ws.Range("C46").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D46").Value = '0.00000000119'
$ws.Range("E46").Value = '  +4.65%  '
$ws.Range("B47").Value = 'RenderToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D47").Value = '1.732'
$ws.Range("E47").Value = '  +2.87%  '
$ws.Range("B48").Value = 'Aptos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D48").Value = '7.055'
$ws.Range("E48").Value = '  -1.50%  '
$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D49").Value = '9.025'
$ws.Range("E49").Value = '  +0.48%  '
$ws.Range("B50").Value = 'Algorand'
$ws.Range("C50").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D50").Value = '0.1154'
$ws.Range("E50").Value = '  +0.89%  '
$ws.Range("D51").Value = '0.3942'
$ws.Range("E51").Value = '  -1.22%  '

# Restore the original (default) cell style for column D so no stray
# number-format / quote-prefix styling is left behind.
$ws.Range("D2:D51").Style = "Normal"
